$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.298.05"
$ws.Range("E2").Value = "  +0.61%  "

$ws.Range("D3").Value = "2.965.76"
$ws.Range("E3").Value = "  +2.39%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'635.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.64%  "

$ws.Range("D6").Value = "'198.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.546"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("E9").Value = "  +2.50%  "

$ws.Range("D10").Value = "2.965.97"
$ws.Range("E10").Value = "  +2.49%  "

$ws.Range("D11").Value = "'0.430"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.87%  "

$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").Value = "'4.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.99%  "

$ws.Range("D14").Value = "3.513.34"
$ws.Range("E14").Value = "  +2.65%  "

$ws.Range("D15").Value = "'28.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.89%  "

$ws.Range("D16").Value = "76.223.09"
$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("D17").Value = "'0.0000186"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").Value = "2.967.45"
$ws.Range("E18").Value = "  +2.59%  "

$ws.Range("D19").Value = "'13.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.78%  "

$ws.Range("D20").Value = "'8.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.08%  "

$ws.Range("D21").Value = "'370.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.55%  "

$ws.Range("D22").Value = "'4.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.07%  "

$ws.Range("D23").Value = "'2.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.05%  "

$ws.Range("D24").Value = "'72.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.31%  "

$ws.Range("D25").Value = "3.122.58"
$ws.Range("E25").Value = "  +2.54%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "'4.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.98%  "

$ws.Range("D28").Value = "'9.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("E29").Value = "  -3.18%  "

$ws.Range("D30").Value = "'0.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("E31").Value = "  +7.08%  "

$ws.Range("D32").Value = "'513.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.26%  "

$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("E34").Value = "  +8.81%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").Value = "'163.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").Value = "'20.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("E38").Value = "  +1.33%  "

$ws.Range("D39").Value = "'0.381"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.94%  "

$ws.Range("D40").Value = "'0.105"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.27%  "

$ws.Range("E41").Value = "  -2.63%  "

$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'181.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.19%  "

$ws.Range("D44").Value = "'42.70"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'4.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.31%  "

$ws.Range("E46").Value = "  -1.92%  "

$ws.Range("D47").Value = "'1.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.42%  "

$ws.Range("E48").Value = "  +8.04%  "

$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("E50").Value = "  -2.43%  "

$ws.Range("D51").Value = "'3.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.16%  "
